# Pinout mit Viktor abgesprochen.
# Adds two new BOM rows (TVS diodes) right after the existing "Trafo" row
# and updates the active selection to the last entered cell (D32).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31: TVS Diode 3.3V SOD323 / 863-ESD7351HT1G
$ws.Range("B31").Value = "TVS Diode 3.3V SOD323"
$ws.Range("D31").Value = "863-ESD7351HT1G "

# Row 32: TVS Diode 5V SOD323 / 833-ESD5V0D3-TP
$ws.Range("B32").Value = "TVS Diode 5V SOD323"
$ws.Range("D32").Value = "833-ESD5V0D3-TP "

# Match the author's final selection/active cell
$ws.Range("D32").Select()
